$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44301
$ws.Range("K2").Value = "Hachiya"
$ws.Range("L2").Value = "Segunda"
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20500
$ws.Range("S2").Value = 1139

# Row 4
$ws.Range("D4").Value = 44342
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 24000
$ws.Range("O4").Value = 25000
$ws.Range("P4").Value = 24500
$ws.Range("S4").Value = 1361

# Row 5
$ws.Range("D5").Value = 44305
$ws.Range("L5").Value = "Segunda"

# Row 6
$ws.Range("D6").Value = 44313
$ws.Range("K6").Value = "Mankaki"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 270
$ws.Range("N6").Value = 21000
$ws.Range("O6").Value = 22000
$ws.Range("P6").Value = 21500
$ws.Range("S6").Value = 1194

# Row 7
$ws.Range("D7").Value = 44699
$ws.Range("L7").Value = "Primera"
$ws.Range("N7").Value = 29000
$ws.Range("O7").Value = 30000
$ws.Range("P7").Value = 29500
$ws.Range("S7").Value = 1639
